# Reorders the "Periodo Mora" values in column E (rows 16-35) from the
# previous descending layout (1802 .. 1607) to the new ascending layout
# (1607 .. 1802), per the updated EC (Estado de Cuenta) database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @(
    "1607", "1608", "1609", "1610", "1611", "1612",
    "1701", "1702", "1703", "1704", "1705", "1706",
    "1707", "1708", "1709", "1710", "1711", "1712",
    "1801", "1802"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}
